$wb = $excel.ActiveWorkbook

# Update violent crime counts to incorporate newly recorded incidents for 2024-04-04.
# Each block below targets one worksheet and writes the corrected cell values
# (crime-category rows and the Total row) for the 2024 column (K) and, where the
# newly-added record shifted a prior-period boundary, the 2023 column (J).

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 1797
$ws.Range("J3").Value = 8076
$ws.Range("K3").Value = 1711
$ws.Range("J4").Value = 1801
$ws.Range("K4").Value = 373
$ws.Range("K6").Value = 2217
$ws.Range("J7").Value = 29269
$ws.Range("K7").Value = 6213

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K4").Value = 3
$ws.Range("K6").Value = 53
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("J4").Value = 96
$ws.Range("J7").Value = 1851

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 130

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 67
$ws.Range("K7").Value = 208

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 42
$ws.Range("K4").Value = 5
$ws.Range("K7").Value = 154

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 39
$ws.Range("K7").Value = 116

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 47
$ws.Range("K5").Value = 10
$ws.Range("K7").Value = 174
$ws.Range("J8").Value = 1851
$ws.Range("K11").Value = 133
$ws.Range("K15").Value = 57
$ws.Range("K18").Value = 46
$ws.Range("K20").Value = 137
$ws.Range("K24").Value = 20
$ws.Range("K29").Value = 301
$ws.Range("K37").Value = 208
$ws.Range("K42").Value = 215
$ws.Range("K43").Value = 61
$ws.Range("K47").Value = 40
$ws.Range("K50").Value = 32
$ws.Range("K51").Value = 72
$ws.Range("K52").Value = 171
$ws.Range("K53").Value = 98
$ws.Range("K54").Value = 100
$ws.Range("J55").Value = 459
$ws.Range("K57").Value = 14
$ws.Range("K64").Value = 42
$ws.Range("K65").Value = 154
$ws.Range("K67").Value = 242
$ws.Range("K71").Value = 16
$ws.Range("K76").Value = 90
$ws.Range("K78").Value = 80
$ws.Range("K79").Value = 166
$ws.Range("K83").Value = 130
$ws.Range("K84").Value = 45
$ws.Range("K85").Value = 310
$ws.Range("K86").Value = 42
$ws.Range("K88").Value = 80
$ws.Range("K90").Value = 54
$ws.Range("K91").Value = 55
$ws.Range("K96").Value = 88
$ws.Range("K99").Value = 116
$ws.Range("J101").Value = 29269
$ws.Range("K101").Value = 6213

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 75
$ws.Range("K7").Value = 242

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range("K2").Value = 16
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 100

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K4").Value = 15
$ws.Range("K6").Value = 97
$ws.Range("K7").Value = 301

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 90

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 57
$ws.Range("K6").Value = 96
$ws.Range("K7").Value = 215

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("J6").Value = 257
$ws.Range("J7").Value = 459

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 20

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 55

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 59
$ws.Range("K3").Value = 58
$ws.Range("K7").Value = 166

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K6").Value = 14
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 40
$ws.Range("K7").Value = 137

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 15
$ws.Range("K7").Value = 46

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 174

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K3").Value = 11
$ws.Range("K4").Value = 8

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 40

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K3").Value = 12
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 32

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K6").Value = 54
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K3").Value = 13
$ws.Range("K7").Value = 47

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("J3").Value = 67
$ws.Range("J6").Value = 166
$ws.Range("K6").Value = 46
$ws.Range("K7").Value = 80

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 10

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("K4").Value = 20
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 18
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 72

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("K2").Value = 5
$ws.Range("K7").Value = 14

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 112
$ws.Range("K3").Value = 101
$ws.Range("K6").Value = 75
$ws.Range("K7").Value = 310

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("K6").Value = 3
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 41
$ws.Range("K6").Value = 79
$ws.Range("K7").Value = 171
